$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = [double]"1.782436333333333"
$ws.Range("H2").Value = [double]"5.347308999999999"
$ws.Range("I2").Value = [double]"0.00914036392049929"
$ws.Range("J2").Value = [double]"0.009140363920499292"
$ws.Range("M2").Value = [double]"227.11144"
$ws.Range("N2").Value = [double]"681.33432"
$ws.Range("O2").Value = [double]"0.8625743548356182"
$ws.Range("P2").Value = [double]"0.8625743548356182"
$ws.Range("Q2").Value = [double]"404.8116823716533"
$ws.Range("R2").Value = [double]"3643.30514134488"
$ws.Range("S2").Value = [double]"0.007884243511687436"
$ws.Range("T2").Value = [double]"0.007884243511687438"

# Row 3
$ws.Range("G3").Value = [double]"1.782436333333333"
$ws.Range("H3").Value = [double]"5.347308999999999"
$ws.Range("I3").Value = [double]"0.00914036392049929"
$ws.Range("J3").Value = [double]"0.009140363920499292"
$ws.Range("O3").Value = [double]"0.001598666154760757"
$ws.Range("P3").Value = [double]"0.001598666154760757"
$ws.Range("Q3").Value = [double]"0.7502642897175555"
$ws.Range("R3").Value = [double]"6.752378607458"
$ws.Range("S3").Value = [double]"1.461239044189856e-05"
$ws.Range("T3").Value = [double]"1.461239044189856e-05"

# Row 4
$ws.Range("G4").Value = [double]"1.782436333333333"
$ws.Range("H4").Value = [double]"5.347308999999999"
$ws.Range("I4").Value = [double]"0.00914036392049929"
$ws.Range("J4").Value = [double]"0.009140363920499292"
$ws.Range("M4").Value = [double]"3.233093"
$ws.Range("N4").Value = [double]"9.699279"
$ws.Range("O4").Value = [double]"0.01227935989749593"
$ws.Range("P4").Value = [double]"0.01227935989749593"
$ws.Range("Q4").Value = [double]"5.762782432245666"
$ws.Range("R4").Value = [double]"51.865041890211"
$ws.Range("S4").Value = [double]"0.0001122378181738977"
$ws.Range("T4").Value = [double]"0.0001122378181738977"

# Row 5
$ws.Range("G5").Value = [double]"1.782436333333333"
$ws.Range("H5").Value = [double]"5.347308999999999"
$ws.Range("I5").Value = [double]"0.00914036392049929"
$ws.Range("J5").Value = [double]"0.009140363920499292"
$ws.Range("M5").Value = [double]"32.52945966666667"
$ws.Range("N5").Value = [double]"97.588379"
$ws.Range("O5").Value = [double]"0.1235476191121251"
$ws.Range("P5").Value = [double]"0.1235476191121251"
$ws.Range("Q5").Value = [double]"57.98169081356788"
$ws.Range("R5").Value = [double]"521.8352173221109"
$ws.Range("S5").Value = [double]"0.001129270200196057"
$ws.Range("T5").Value = [double]"0.001129270200196057"

# Row 6
$ws.Range("H6").Value = [double]"564.692825"
$ws.Range("I6").Value = [double]"0.965251479537618"
$ws.Range("J6").Value = [double]"0.965251479537618"
$ws.Range("M6").Value = [double]"227.11144"
$ws.Range("N6").Value = [double]"681.33432"
$ws.Range("O6").Value = [double]"0.8625743548356182"
$ws.Range("P6").Value = [double]"0.8625743548356182"
$ws.Range("Q6").Value = [double]"42749.40021447267"
$ws.Range("R6").Value = [double]"384744.601930254"
$ws.Range("S6").Value = [double]"0.8326011722162868"
$ws.Range("T6").Value = [double]"0.8326011722162868"

# Row 7
$ws.Range("H7").Value = [double]"564.692825"
$ws.Range("I7").Value = [double]"0.965251479537618"
$ws.Range("J7").Value = [double]"0.965251479537618"
$ws.Range("O7").Value = [double]"0.001598666154760757"
$ws.Range("P7").Value = [double]"0.001598666154760757"
$ws.Range("R7").Value = [double]"713.0726410826501"
$ws.Range("S7").Value = [double]"0.001543114871169535"
$ws.Range("T7").Value = [double]"0.001543114871169535"

# Row 8
$ws.Range("H8").Value = [double]"564.692825"
$ws.Range("I8").Value = [double]"0.965251479537618"
$ws.Range("J8").Value = [double]"0.965251479537618"
$ws.Range("M8").Value = [double]"3.233093"
$ws.Range("N8").Value = [double]"9.699279"
$ws.Range("O8").Value = [double]"0.01227935989749593"
$ws.Range("P8").Value = [double]"0.01227935989749593"
$ws.Range("Q8").Value = [double]"608.5681398859084"
$ws.Range("R8").Value = [double]"5477.113258973175"
$ws.Range("S8").Value = [double]"0.01185267030883284"
$ws.Range("T8").Value = [double]"0.01185267030883284"

# Row 9
$ws.Range("H9").Value = [double]"564.692825"
$ws.Range("I9").Value = [double]"0.965251479537618"
$ws.Range("J9").Value = [double]"0.965251479537618"
$ws.Range("M9").Value = [double]"32.52945966666667"
$ws.Range("N9").Value = [double]"97.588379"
$ws.Range("O9").Value = [double]"0.1235476191121251"
$ws.Range("P9").Value = [double]"0.1235476191121251"
$ws.Range("Q9").Value = [double]"6123.05082496452"
$ws.Range("R9").Value = [double]"55107.45742468067"
$ws.Range("S9").Value = [double]"0.1192545221413289"
$ws.Range("T9").Value = [double]"0.1192545221413289"

# Row 10
$ws.Range("G10").Value = [double]"4.870778333333333"
$ws.Range("H10").Value = [double]"14.612335"
$ws.Range("I10").Value = [double]"0.02497743437460768"
$ws.Range("J10").Value = [double]"0.02497743437460768"
$ws.Range("M10").Value = [double]"227.11144"
$ws.Range("N10").Value = [double]"681.33432"
$ws.Range("O10").Value = [double]"0.8625743548356182"
$ws.Range("P10").Value = [double]"0.8625743548356182"
$ws.Range("Q10").Value = [double]"1106.209481204133"
$ws.Range("R10").Value = [double]"9955.8853308372"
$ws.Range("S10").Value = [double]"0.02154489434112621"
$ws.Range("T10").Value = [double]"0.02154489434112621"

# Row 11
$ws.Range("G11").Value = [double]"4.870778333333333"
$ws.Range("H11").Value = [double]"14.612335"
$ws.Range("I11").Value = [double]"0.02497743437460768"
$ws.Range("J11").Value = [double]"0.02497743437460768"
$ws.Range("O11").Value = [double]"0.001598666154760757"
$ws.Range("P11").Value = [double]"0.001598666154760757"
$ws.Range("Q11").Value = [double]"2.050211263252222"
$ws.Range("R11").Value = [double]"18.45190136927"
$ws.Range("S11").Value = [double]"3.993057896744321e-05"
$ws.Range("T11").Value = [double]"3.99305789674432e-05"

# Row 12
$ws.Range("G12").Value = [double]"4.870778333333333"
$ws.Range("H12").Value = [double]"14.612335"
$ws.Range("I12").Value = [double]"0.02497743437460768"
$ws.Range("J12").Value = [double]"0.02497743437460768"
$ws.Range("M12").Value = [double]"3.233093"
$ws.Range("N12").Value = [double]"9.699279"
$ws.Range("O12").Value = [double]"0.01227935989749593"
$ws.Range("P12").Value = [double]"0.01227935989749593"
$ws.Range("Q12").Value = [double]"15.74767933405167"
$ws.Range("R12").Value = [double]"141.729114006465"
$ws.Range("S12").Value = [double]"0.0003067069060018939"
$ws.Range("T12").Value = [double]"0.0003067069060018939"

# Row 13
$ws.Range("G13").Value = [double]"4.870778333333333"
$ws.Range("H13").Value = [double]"14.612335"
$ws.Range("I13").Value = [double]"0.02497743437460768"
$ws.Range("J13").Value = [double]"0.02497743437460768"
$ws.Range("M13").Value = [double]"32.52945966666667"
$ws.Range("N13").Value = [double]"97.588379"
$ws.Range("O13").Value = [double]"0.1235476191121251"
$ws.Range("P13").Value = [double]"0.1235476191121251"
$ws.Range("Q13").Value = [double]"158.4437873394405"
$ws.Range("R13").Value = [double]"1425.994086054965"
$ws.Range("S13").Value = [double]"0.003085902548512131"
$ws.Range("T13").Value = [double]"0.003085902548512131"

# Row 14
$ws.Range("E14").Value = [double]"3.0"
$ws.Range("F14").Value = [double]"1.0"
$ws.Range("G14").Value = [double]"0.1229953333333333"
$ws.Range("H14").Value = [double]"0.368986"
$ws.Range("I14").Value = [double]"0.0006307221672750447"
$ws.Range("J14").Value = [double]"0.0006307221672750447"
$ws.Range("M14").Value = [double]"227.11144"
$ws.Range("N14").Value = [double]"681.33432"
$ws.Range("O14").Value = [double]"0.8625743548356182"
$ws.Range("P14").Value = [double]"0.8625743548356182"
$ws.Range("Q14").Value = [double]"27.93364726661334"
$ws.Range("R14").Value = [double]"251.4028253995201"
$ws.Range("S14").Value = [double]"0.0005440447665177946"
$ws.Range("T14").Value = [double]"0.0005440447665177946"

# Row 15
$ws.Range("E15").Value = [double]"3.0"
$ws.Range("F15").Value = [double]"1.0"
$ws.Range("G15").Value = [double]"0.1229953333333333"
$ws.Range("H15").Value = [double]"0.368986"
$ws.Range("I15").Value = [double]"0.0006307221672750447"
$ws.Range("J15").Value = [double]"0.0006307221672750447"
$ws.Range("O15").Value = [double]"0.001598666154760757"
$ws.Range("P15").Value = [double]"0.001598666154760757"
$ws.Range("Q15").Value = [double]"0.05177127770355557"
$ws.Range("R15").Value = [double]"0.4659414993320001"
$ws.Range("S15").Value = [double]"1.008314181879967e-06"
$ws.Range("T15").Value = [double]"1.008314181879966e-06"

# Row 16
$ws.Range("E16").Value = [double]"3.0"
$ws.Range("F16").Value = [double]"1.0"
$ws.Range("G16").Value = [double]"0.1229953333333333"
$ws.Range("H16").Value = [double]"0.368986"
$ws.Range("I16").Value = [double]"0.0006307221672750447"
$ws.Range("J16").Value = [double]"0.0006307221672750447"
$ws.Range("M16").Value = [double]"3.233093"
$ws.Range("N16").Value = [double]"9.699279"
$ws.Range("O16").Value = [double]"0.01227935989749593"
$ws.Range("P16").Value = [double]"0.01227935989749593"
$ws.Range("Q16").Value = [double]"0.3976553512326667"
$ws.Range("R16").Value = [double]"3.578898161094"
$ws.Range("S16").Value = [double]"7.744864487298906e-06"
$ws.Range("T16").Value = [double]"7.744864487298904e-06"

# Row 17
$ws.Range("E17").Value = [double]"3.0"
$ws.Range("F17").Value = [double]"1.0"
$ws.Range("G17").Value = [double]"0.1229953333333333"
$ws.Range("H17").Value = [double]"0.368986"
$ws.Range("I17").Value = [double]"0.0006307221672750447"
$ws.Range("J17").Value = [double]"0.0006307221672750447"
$ws.Range("M17").Value = [double]"32.52945966666667"
$ws.Range("N17").Value = [double]"97.588379"
$ws.Range("O17").Value = [double]"0.1235476191121251"
$ws.Range("P17").Value = [double]"0.1235476191121251"
$ws.Range("Q17").Value = [double]"4.00097173485489"
$ws.Range("R17").Value = [double]"36.008745613694"
$ws.Range("S17").Value = [double]"7.79242220880713e-05"
$ws.Range("T17").Value = [double]"7.79242220880713e-05"

Write-Output "applied changes"